$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 736, shifting rows 736:777 down to 737:778.
$ws.Rows.Item(736).Insert()

# Populate the newly inserted row 736 with the new log entry
# (2026/01/29, 木, 22:00, ranking 201).
#
# Column A stores dates as plain text ("2026/01/29", no special number
# format). Assigning that string straight to .Value would get parsed as
# a date serial instead of staying text. Row 735 already holds the exact
# same date text we need here (it's one more reading on 2026/01/29,
# right before the log resumes in 2026/12), so copy that cell (value +
# format) down into the new row instead of typing the value in - this
# keeps the cell a plain text cell with the default style, matching the
# rest of the column.
$ws.Cells.Item(735, 1).Copy($ws.Cells.Item(736, 1))

$ws.Cells.Item(736, 2).Value = "木"
$ws.Cells.Item(736, 3).Value = 22
$ws.Cells.Item(736, 4).Value = 201
